$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Delete row 49 first (the obsolete "Value from main.c: <value>" entry).
# Everything below shifts up by one row.
$ws.Rows("49:49").Delete()

# After that shift, the two trailing duplicate rows (originally 74 and 75)
# now sit at 73 and 74; remove them too.
$ws.Rows("73:74").Delete()
